$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 44860
$ws.Range("J3").Value = 1100
$ws.Range("K3").Value = 1500
$ws.Range("L3").Value = 1700
$ws.Range("M3").Value = 1609
$ws.Range("P3").Value = 1609
# Row 4
$ws.Range("D4").Value = 44519
$ws.Range("J4").Value = 250
$ws.Range("K4").Value = 1200
$ws.Range("L4").Value = 1300
$ws.Range("M4").Value = 1240
$ws.Range("P4").Value = 1240
# Row 5
$ws.Range("D5").Value = 44510
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 600
$ws.Range("K5").Value = 1300
$ws.Range("L5").Value = 1400
$ws.Range("M5").Value = 1350
$ws.Range("P5").Value = 1350
# Row 6
$ws.Range("D6").Value = 44496
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("J6").Value = 550
$ws.Range("K6").Value = 1500
$ws.Range("M6").Value = 1773
$ws.Range("N6").Value = '$/paquete'
$ws.Range("P6").Value = 1773
# Row 7
$ws.Range("D7").Value = 44876
$ws.Range("J7").Value = 350
$ws.Range("L7").Value = 1600
$ws.Range("M7").Value = 1557
$ws.Range("N7").Value = '$/kilo'
$ws.Range("P7").Value = 1557
# Row 8
$ws.Range("D8").Value = 44489
$ws.Range("K8").Value = 1400
$ws.Range("L8").Value = 1500
$ws.Range("M8").Value = 1450
$ws.Range("P8").Value = 1450
# Row 10
$ws.Range("D10").Value = 44881
$ws.Range("J10").Value = 200
$ws.Range("K10").Value = 2600
$ws.Range("L10").Value = 2700
$ws.Range("M10").Value = 2650
$ws.Range("P10").Value = 2650
# Row 11
$ws.Range("D11").Value = 44881
$ws.Range("I11").Value = "Segunda"
$ws.Range("J11").Value = 100
$ws.Range("K11").Value = 2400
$ws.Range("L11").Value = 2400
$ws.Range("M11").Value = 2400
$ws.Range("P11").Value = 2400
# Row 12
$ws.Range("D12").Value = 44875
$ws.Range("J12").Value = 300
$ws.Range("L12").Value = 1600
$ws.Range("M12").Value = 1550
$ws.Range("P12").Value = 1550
# Row 13
$ws.Range("D13").Value = 44868
$ws.Range("J13").Value = 1000
$ws.Range("K13").Value = 1200
$ws.Range("L13").Value = 1300
$ws.Range("M13").Value = 1250
$ws.Range("O13").Value = "Región del Maule"
$ws.Range("P13").Value = 1250
# Row 14
$ws.Range("D14").Value = 44868
$ws.Range("I14").Value = "Segunda"
$ws.Range("J14").Value = 200
$ws.Range("K14").Value = 1000
$ws.Range("L14").Value = 1000
$ws.Range("M14").Value = 1000
$ws.Range("O14").Value = "Región del Maule"
$ws.Range("P14").Value = 1000
# Row 15
$ws.Range("D15").Value = 44511
$ws.Range("J15").Value = 600
$ws.Range("K15").Value = 1300
$ws.Range("L15").Value = 1400
$ws.Range("M15").Value = 1350
$ws.Range("P15").Value = 1350
# Row 16
$ws.Range("D16").Value = 44839
$ws.Range("J16").Value = 500
$ws.Range("K16").Value = 1700
$ws.Range("L16").Value = 1800
$ws.Range("M16").Value = 1760
$ws.Range("P16").Value = 1760
# Row 17
$ws.Range("D17").Value = 44468
$ws.Range("H17").Value = "Verde"
$ws.Range("J17").Value = 500
$ws.Range("K17").Value = 1800
$ws.Range("L17").Value = 2000
$ws.Range("M17").Value = 1920
$ws.Range("O17").Value = "Provincia de Linares"
$ws.Range("P17").Value = 1920
# Row 18
$ws.Range("D18").Value = 44477
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 500
$ws.Range("K18").Value = 1400
$ws.Range("L18").Value = 1500
$ws.Range("M18").Value = 1460
$ws.Range("O18").Value = "Provincia de Linares"
$ws.Range("P18").Value = 1460
# Row 19
$ws.Range("D19").Value = 44545
$ws.Range("J19").Value = 550
$ws.Range("K19").Value = 1700
$ws.Range("L19").Value = 1800
$ws.Range("M19").Value = 1755
$ws.Range("P19").Value = 1755
